$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("M2").Value = 29.47402433333333
$ws.Range("N2").Value = 88.422073
$ws.Range("O2").Value = 0.295877356230023
$ws.Range("P2").Value = 0.295877356230023
$ws.Range("Q2").Value = 9.230645466688999
$ws.Range("R2").Value = 83.07580920020099
$ws.Range("S2").Value = 0.008010774669212066
$ws.Range("T2").Value = 0.008010774669212068
$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("O3").Value = 0.1818061388681701
$ws.Range("P3").Value = 0.1818061388681701
$ws.Range("Q3").Value = 5.671904173211
$ws.Range("R3").Value = 51.047137558899
$ws.Range("S3").Value = 0.004922336844256976
$ws.Range("T3").Value = 0.004922336844256977
$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 7.238098333333333
$ws.Range("N4").Value = 21.714295
$ws.Range("O4").Value = 0.07266023040422054
$ws.Range("P4").Value = 0.07266023040422054
$ws.Range("Q4").Value = 2.266820397935
$ws.Range("R4").Value = 20.401383581415
$ws.Range("S4").Value = 0.001967250013984611
$ws.Range("T4").Value = 0.001967250013984611
$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("M5").Value = 44.79281599999999
$ws.Range("N5").Value = 134.378448
$ws.Range("O5").Value = 0.4496562744975863
$ws.Range("P5").Value = 0.4496562744975864
$ws.Range("Q5").Value = 14.028169322064
$ws.Range("R5").Value = 126.253523898576
$ws.Range("S5").Value = 0.01217428443830344
$ws.Range("T5").Value = 0.01217428443830344
$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("M6").Value = 29.47402433333333
$ws.Range("N6").Value = 88.422073
$ws.Range("O6").Value = 0.295877356230023
$ws.Range("P6").Value = 0.295877356230023
$ws.Range("Q6").Value = 238.0670391634125
$ws.Range("R6").Value = 2142.603352470712
$ws.Range("S6").Value = 0.2066054225337561
$ws.Range("T6").Value = 0.2066054225337562
$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("O7").Value = 0.1818061388681701
$ws.Range("P7").Value = 0.1818061388681701
$ws.Range("S7").Value = 0.1269517026199437
$ws.Range("T7").Value = 0.1269517026199437
$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 7.238098333333333
$ws.Range("N8").Value = 21.714295
$ws.Range("O8").Value = 0.07266023040422054
$ws.Range("P8").Value = 0.07266023040422054
$ws.Range("Q8").Value = 58.46343274683112
$ws.Range("R8").Value = 526.17089472148
$ws.Range("S8").Value = 0.05073723043676694
$ws.Range("T8").Value = 0.05073723043676695
$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("M9").Value = 44.79281599999999
$ws.Range("N9").Value = 134.378448
$ws.Range("O9").Value = 0.4496562744975863
$ws.Range("P9").Value = 0.4496562744975864
$ws.Range("Q9").Value = 361.7996972626347
$ws.Range("R9").Value = 3256.197275363712
$ws.Range("S9").Value = 0.3139862602912553
$ws.Range("T9").Value = 0.3139862602912554
$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("M10").Value = 29.47402433333333
$ws.Range("N10").Value = 88.422073
$ws.Range("O10").Value = 0.295877356230023
$ws.Range("P10").Value = 0.295877356230023
$ws.Range("Q10").Value = 85.40822629114456
$ws.Range("R10").Value = 768.674036620301
$ws.Range("S10").Value = 0.07412114983556486
$ws.Range("T10").Value = 0.07412114983556486
$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("O11").Value = 0.1818061388681701
$ws.Range("P11").Value = 0.1818061388681701
$ws.Range("Q11").Value = 52.48032511653323
$ws.Range("R11").Value = 472.322926048799
$ws.Range("S11").Value = 0.04554481705452578
$ws.Range("T11").Value = 0.04554481705452578
$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 7.238098333333333
$ws.Range("N12").Value = 21.714295
$ws.Range("O12").Value = 0.07266023040422054
$ws.Range("P12").Value = 0.07266023040422054
$ws.Range("Q12").Value = 20.97416808032389
$ws.Range("R12").Value = 188.767512722915
$ws.Range("S12").Value = 0.01820233860914635
$ws.Range("T12").Value = 0.01820233860914635
$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("M13").Value = 44.79281599999999
$ws.Range("N13").Value = 134.378448
$ws.Range("O13").Value = 0.4496562744975863
$ws.Range("P13").Value = 0.4496562744975864
$ws.Range("Q13").Value = 129.7981884617973
$ws.Range("R13").Value = 1168.183696156176
$ws.Range("S13").Value = 0.1126447813418564
$ws.Range("T13").Value = 0.1126447813418564
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("M14").Value = 29.47402433333333
$ws.Range("N14").Value = 88.422073
$ws.Range("O14").Value = 0.295877356230023
$ws.Range("P14").Value = 0.295877356230023
$ws.Range("Q14").Value = 8.227280905658889
$ws.Range("R14").Value = 74.04552815093
$ws.Range("S14").Value = 0.007140009191489931
$ws.Range("T14").Value = 0.007140009191489933
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("O15").Value = 0.1818061388681701
$ws.Range("P15").Value = 0.1818061388681701
$ws.Range("Q15").Value = 5.055372245785556
$ws.Range("R15").Value = 45.49835021207
$ws.Range("S15").Value = 0.004387282349443645
$ws.Range("T15").Value = 0.004387282349443646
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 7.238098333333333
$ws.Range("N16").Value = 21.714295
$ws.Range("O16").Value = 0.07266023040422054
$ws.Range("P16").Value = 0.07266023040422054
$ws.Range("Q16").Value = 2.020418641772222
$ws.Range("R16").Value = 18.18376777595
$ws.Range("S16").Value = 0.001753411344322632
$ws.Range("T16").Value = 0.001753411344322632
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("M17").Value = 44.79281599999999
$ws.Range("N17").Value = 134.378448
$ws.Range("O17").Value = 0.4496562744975863
$ws.Range("P17").Value = 0.4496562744975864
$ws.Range("Q17").Value = 12.50331734885333
$ws.Range("R17").Value = 112.52985613968
$ws.Range("S17").Value = 0.01085094842617128
$ws.Range("T17").Value = 0.01085094842617128
